# "menage + class + recalcul"
# Réorganisation de la matrice temps : 5 équipes (A-E) de 4 joueurs
# -> 3 équipes (A-C) de 7/8 joueurs, colonnes F:G supprimées,
# ajout d'un nouveau joueur (D. Wade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repartir sur une base propre : on efface toute la zone utilisée actuelle
$ws.Range("A1:G5").Delete()

# En-tete
$ws.Range("A1").Value = "Equipes"

$ws.Range("C1").Value = "A"
$ws.Range("D1").Value = "B"
$ws.Range("E1").Value = "C"

$ws.Range("C1:E1").Font.Bold = $true

# Equipe A (colonne C) : anciennes equipes A + B fusionnees
$ws.Range("C2").Value = "F. Mayweather"
$ws.Range("C3").Value = "C. Ronaldo"
$ws.Range("C4").Value = "L. Messi"
$ws.Range("C5").Value = "K. Bryant"
$ws.Range("C6").Value = "R. Federer"
$ws.Range("C7").Value = "P. Mickelson"
$ws.Range("C8").Value = "R. Nadal"

# Equipe B (colonne D) : anciennes equipes C + D fusionnees
$ws.Range("D2").Value = "M. Ryan"
$ws.Range("D3").Value = "M. Pacquiao"
$ws.Range("D4").Value = "Z. Ibrahimović"
$ws.Range("D5").Value = "D. Rose"
$ws.Range("D6").Value = "G. Bale"
$ws.Range("D7").Value = "R. Falcao"
$ws.Range("D8").Value = "M. Özil"

# Equipe C (colonne E) : ancienne equipe E + nouveau joueur
$ws.Range("E2").Value = "N. Djokovic"
$ws.Range("E3").Value = "M. Stafford"
$ws.Range("E4").Value = "L. Hamilton"
$ws.Range("E5").Value = "K. Durant"
$ws.Range("E6").Value = "F. Alonso"
$ws.Range("E7").Value = "M. Singh"
$ws.Range("E8").Value = "D. Wade"

$wb.Application.CalculateFull()
